# Update the distance codes (D80->D86, D51->D55, D64->D69) and the
# "S30" size code (->S31) throughout the order sheet. These tokens show
# up embedded inside several text columns (Condition, Filename_Left,
# Filename_Right, Distance, Size), so every cell in those columns needs
# a substring replacement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

# Map header name -> column index by reading row 1.
$targetHeaders = @("Condition", "Filename_Left", "Filename_Right", "Distance", "Size")
$targetCols = @()
for ($c = 1; $c -le $colCount; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($targetHeaders -contains $header) {
        $targetCols += $c
    }
}

for ($r = 2; $r -le $rowCount; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $newVal = $val.Replace("D80", "D86").Replace("D51", "D55").Replace("D64", "D69").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
